$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11 for "One or more PowerBI reports" — this
# shifts the old rows 11-14 down to 12-15 and auto-adjusts the SUM() ranges
# in the totals row.
$ws.Rows(11).Insert()

# --- Fill in the new row 11 (PowerBI reports line) ---
$ws.Range("B11").Value = 10
$ws.Range("D11").Value = "One or more PowerBI reports"
$ws.Range("E11").Value = 1

# --- Update "Max Points" (column B) figures throughout the table ---
$ws.Range("B3").Value = 50
$ws.Range("B4").Value = 10
$ws.Range("B5").Value = 50
$ws.Range("B7").Value = 10
$ws.Range("B8").Value = 10
$ws.Range("B9").Value = 10
$ws.Range("B10").Value = 10
$ws.Range("B12").Value = 50
$ws.Range("B13").Value = 40

# --- Fix capitalisation of two component labels ---
$ws.Range("D9").Value = "One or more Excel reports"
$ws.Range("D10").Value = "One or more Report Builder reports"

# --- Keep the selection where the author left it ---
$ws.Range("B14").Select()
